# "add start/resume on env and fix pagerank db"
# Adds a small "Jumlah Request" / "Rata-rata Response Time" load-test
# summary table in columns L:M (rows 7-10) next to the existing API
# response-time table, matching the formatting of the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data -------------------------------------------------------
# Enter the response-time values first, then the header row, then the
# request-count column, so new shared-string entries land in the same
# order as the source edit.
$ws.Range("M8").Value = "581 ms"
$ws.Range("M9").Value = "796 ms"
$ws.Range("M10").Value = "1826 ms"
$ws.Range("M7").Value = "Rata-rata Response Time"
$ws.Range("L7").Value = "Jumlah Request"
$ws.Range("L8").Value = 100
$ws.Range("L9").Value = 1000
$ws.Range("L10").Value = 10000

# --- Formatting -------------------------------------------------------
# Reuse the formatting already used in the adjacent table instead of
# building new styles from scratch: the bold header cell (I8) for the
# L7:M7 header row, and the plain bordered cell (I9) for the data rows.
$ws.Range("I8").Copy()
$ws.Range("L7:M7").PasteSpecial(-4122)
$ws.Range("I9").Copy()
$ws.Range("L8:M10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths for the two new columns.
$ws.Columns("L").ColumnWidth = 16.5
$ws.Columns("M").ColumnWidth = 23.5

# Row 10 relaxes back to the default single-line height.
$ws.Rows("10:10").RowHeight = 15.75

# Final selection, matching the end of the editing session.
[void]$ws.Range("L11").Select()
